$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of progress-report data (row 12), matching the formatting already
# used by the identical "Lambda Functions" rows above (C/D columns centered /
# left-aligned respectively, F/G as start/end time-of-day values).
$ws.Range("C12").HorizontalAlignment = -4108
$ws.Range("C12").Value = 6

$ws.Range("D12").HorizontalAlignment = -4131
$ws.Range("D12").Value = "Lambda Functions"

$ws.Range("E12").Value = "permissions, Policies, Eventand Context Objects, Destinations,logging & Monitoring,X-Ray"

$ws.Range("F12:G12").NumberFormat = "h:mm"
$ws.Range("F12").Value = 0.64236111111111116
$ws.Range("G12").Value = 0.69444444444444442

# Update the saved selection to match where the author left off editing
$ws.Range("E14").Select()
